$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ayats = "Surah Baqarah, 145 - 160"
$tags = "Patience in life, Elderly philosophy, Wisdom vs Knowledge"
$content = @'
h1: How hard is patience? 
p: Really hard. Initially it is really hard. Then it dulls you down. When it breaks you down and all your pride is dusted away, it shows its colors. You do not react to losses, fights and humiliations. Patience gives you ease in hardships.
h3: Baba ji (Elderly) Philosophy
p: We all have those aged and humbled elders in our homes. Can’t even walk properly. He wakes up early in the morning. Takes shower, wears scarf, gown and warm socks. Finds his crutches and slowly walks to the mosque. It is winter / summer does not matter. He just does everything he has been doing for so long, everyday. Yet, as he reaches the mosque, he is in time. 
p: If patience is so hard, why does a person on crutches prefer a hard life over his warm bed’s comfort in chilling winter? May be because patience is not as hard as it looks. May be it is the end product of our Eman (Faith).
h3: How easy is patience?
p: If it is do able, then how do able is it?
quote: And We will surely test you with something of <u>fear</u> and <u>hunger</u> and a <u>loss of wealth</u> and<u> lives</u> and<u> fruits</u>, but give good tidings to the patient. <br>- Surah Baqarah verse 155
p: Things really dearer to us will be taken from us. Some of us have already lost them. Lost our car in accident, parents in early age, freedom or say poverty struck us. If your Blood pressure went Up! or it was misery and hopelessness you felt, then it is troubling. Allah (swt) is giving good news to people who have pateince. People who sat back. Smiled and embraced the calamity upon them. Forgot it and resumed their life to their goals while claiming:-
quote: Indeed we belong to Allah, and indeed to Him we will return. <br> - Surah Baqarah verse 156
p: .. are successful. What do they get for showing resilience to the loss?
quote: Those are the ones upon whom are blessings from their Lord and mercy. And it is those who are the [rightly] guided. <br>- Surah Baqarah verse 157
p: May be patience is a <b>skill</b> that needs mastery. May be it is a prerequisite to a happy life. Like, anyother skill in this world, this trait also needs learning, class discussions, making groups, doing assignments, undergoing small exams here and there. Like a skill is taught in some institute. It is just another skill and it can be learnt. <a href="https://www.huffpost.com/entry/patience-tips_n_5843928">Click here to learn 5 Tricks To Becoming A More Patient Person.</a>
quote: Strange are the ways of a believer for there is good in every affair of his, and this is not the case with anyone else except in the case of a believer; for if he has an occasion to feel delight, he thanks (God); thus there is a good for him in it, and if he gets into trouble and shows resignation (and endures it patiently), there is a good for him in it. <br> - Sahih Muslim, 2999
h3: Rocket science?
p: Heard it every where. Since childhood reading in our books. 
p.b-left: Say 5 prayers a day. 
p.b-left: Do not tell a lie. 
p.b-left: Be calm, be happy and be hopeful. 
p.b-left: Expect good and have good. 
p: We are consuming so many books, but not acting upon anything we are reading. Book reading has become a race <b>'who reads more books faster'</b>.  More in information, but low in wisdom. We are running with our eyes closed to a dead end. A bigger house, a car, a secure future. Cramped current affairs with their conclusions, references, names of anchor persons, time, location, their affect on future. Our head is buzzing with loud un-necessary noises. We are in a state of a long never ending race.
p: May be we need to stop running and think for a moment where are we heading. May be read less books. Sit calmly with our families. Resort to more patience. <span class="lavendar">Act upon what we read in our 5th grade </span>and stop running on the wrong track.   
quote: Wisdom is the right application of knowledge; and true education...is the application of knowledge to the development of a noble and Godlike character. <br> - David O. McKay, Gospel Ideals
h3: Conclusion
p: Calamities are going to come. Patience or no patience. Everyone of us will be tested. Winners are people, who faced the calamities, stayed patient and became Quaid-e-Azam, Gandhi, Edhi and every brave successful person around us. 
'@

# Day 15 row: Ser, Date, Ayats, Content, Author, Tags
$ws.Range("A16").Value = 15
$ws.Range("B16").NumberFormat = $ws.Range("B15").NumberFormat
$ws.Range("B16").Value = 43845
$ws.Range("C16").Value = $ayats
$ws.Range("F16").Value = $tags
$ws.Range("D16").Value = $content
$ws.Range("E16").Value = "Qasim Ali"

$ws.Rows.Item(16).RowHeight = 409.6

$null = $ws.Range("D16").Select()
